$d = $word.ActiveDocument

# Insert four new paragraphs, each containing "Passw0rd@123", before the
# document's original first paragraph (which already reads "Passw0rd@123").
for ($i = 0; $i -lt 4; $i++) {
    $d.Paragraphs(1).Range.InsertBefore("Passw0rd@123`r`n")
}
